$wb = $excel.ActiveWorkbook

# --- Sheet "CPPbI": add a new "process emissions" column (C), mirroring column B ---
$wsCPPbI = $wb.Worksheets.Item("CPPbI")

$wsCPPbI.Range("B1").Value2 = "energy related emissions"
$wsCPPbI.Range("C1").Value2 = "process emissions"

$wsCPPbI.Range("C2").Value2 = $wsCPPbI.Range("B2").Value2
$wsCPPbI.Range("C3").Value2 = $wsCPPbI.Range("B3").Value2
$wsCPPbI.Range("C4").Value2 = $wsCPPbI.Range("B4").Value2
$wsCPPbI.Range("C5").Value2 = $wsCPPbI.Range("B5").Value2
$wsCPPbI.Range("C6").Value2 = $wsCPPbI.Range("B6").Value2
$wsCPPbI.Range("C7").Value2 = $wsCPPbI.Range("B7").Value2
$wsCPPbI.Range("C8").Value2 = $wsCPPbI.Range("B8").Value2
$wsCPPbI.Range("C9").Value2 = $wsCPPbI.Range("B9").Value2

# --- Sheet "About": add Oregon label in B1 ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("B1").Value2 = "Oregon"

$wb.Save()
